$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-10 Saturday", "2024-08-11 Sunday"),
    @("91÷3=", "74÷4="),
    @("40÷4=", "95÷5="),
    @("87÷8=", "84÷8="),
    @("94÷9=", "90÷8="),
    @("17÷4=", "50÷9="),
    @("93÷9=", "12÷5="),
    @("66÷8=", "83÷9="),
    @("28÷5=", "12÷8="),
    @("52÷2=", "19÷2="),
    @("51÷7=", "50÷4="),
    @("85÷4=", "32÷4="),
    @("61÷9=", "83÷9="),
    @("10÷6=", "84÷7="),
    @("86÷8=", "97÷5="),
    @("86÷7=", "68÷8="),
    @("19÷8=", "82÷3="),
    @("14÷4=", "77÷5="),
    @("87÷2=", "34÷4="),
    @("70÷4=", "54÷2="),
    @("80÷6=", "30÷6="),
    @("73÷4=", "85÷8="),
    @("90÷5=", "96÷7="),
    @("11÷6=", "28÷7="),
    @("10÷3=", "32÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
